# Update "想去人数" (interested-count) figures in the "展览" and "全部类型"
# worksheets to reflect newer scrape totals.

$wb = $excel.ActiveWorkbook

# Map of row -> new F-column value, applied identically to both sheets.
$updates = @{
    2  = 3116
    4  = 2872
    8  = 1536
    14 = 402
    21 = 2866
    22 = 345
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

$wb.Save()
